$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.943.06"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -0.14%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.636.75"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -0.67%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'212.48"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -0.60%  "
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'  -0.17%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'23.33"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -1.31%  "
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'  -2.20%  "
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'  -0.21%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.0882"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +1.16%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'1.869.20"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -0.64%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'1.630.77"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -1.08%  "
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'  -0.22%  "
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'  +0.93%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'65.31"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -0.62%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'27.949.39"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  -0.10%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'231.14"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -0.32%  "
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'  -0.43%  "
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'  -1.76%  "
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'  -0.11%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'10.42"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -2.71%  "
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'  -0.51%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'2.08"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -3.73%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'153.63"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +1.04%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'6.97"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +0.74%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'  -0.40%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  -0.67%  "
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = "'  -0.83%  "
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "'  -0.66%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'  +1.51%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'1.408.03"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -3.17%  "
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "'  -1.42%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "'  +1.38%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'  +1.61%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'0.0170"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +0.38%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.563"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +0.21%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.930"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +1.49%  "
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'  -1.49%  "
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'  +0.30%  "
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'67.16"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -3.25%  "
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'  +2.54%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'1.82"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  +2.22%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'  -1.52%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'1.778.57"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -0.69%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'88.00"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -1.11%  "
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'  -0.38%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  -0.36%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'7.56"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -2.47%  "
$c.Style = "Normal"
